# Update "想去人数" (number of people interested) values in column F
# across the "展览", "演出" and "全部类型" sheets, reflecting the
# latest scrape output (gh-pages update at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value  = 1671
$wsExpo.Range("F3").Value  = 9101
$wsExpo.Range("F4").Value  = 112
$wsExpo.Range("F6").Value  = 706
$wsExpo.Range("F7").Value  = 1371
$wsExpo.Range("F9").Value  = 57
$wsExpo.Range("F10").Value = 94
$wsExpo.Range("F11").Value = 5902
$wsExpo.Range("F13").Value = 387
$wsExpo.Range("F15").Value = 4429
$wsExpo.Range("F18").Value = 1148
$wsExpo.Range("F19").Value = 28
$wsExpo.Range("F21").Value = 22
$wsExpo.Range("F22").Value = 256
$wsExpo.Range("F23").Value = 14
$wsExpo.Range("F24").Value = 2745
$wsExpo.Range("F25").Value = 126

# --- Sheet "演出" (Performances) ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 30

# --- Sheet "全部类型" (All types, combined listing) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value  = 1671
$wsAll.Range("F3").Value  = 9101
$wsAll.Range("F4").Value  = 112
$wsAll.Range("F5").Value  = 30
$wsAll.Range("F7").Value  = 706
$wsAll.Range("F8").Value  = 1371
$wsAll.Range("F10").Value = 57
$wsAll.Range("F11").Value = 94
$wsAll.Range("F12").Value = 5902
$wsAll.Range("F14").Value = 387
$wsAll.Range("F16").Value = 4429
$wsAll.Range("F19").Value = 1148
$wsAll.Range("F20").Value = 28
$wsAll.Range("F22").Value = 22
$wsAll.Range("F23").Value = 256
$wsAll.Range("F24").Value = 14
$wsAll.Range("F25").Value = 2745
$wsAll.Range("F27").Value = 126
